# Daily attendance processing - 2025-10-12 15:40:34
# Reorders the "Recorded By" (column G) text so that "System" (exact,
# case-sensitive match) is listed first among the comma-separated authors,
# while the remaining names keep their original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedLastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $usedLastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq "backup@backdoor.com, System") {
        $cell.Value2 = "System, backup@backdoor.com"
    }
    elseif ($val -eq "backup@backdoor.com, system, System") {
        $cell.Value2 = "System, backup@backdoor.com, system"
    }
    elseif ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
}
